$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row => [newD (or $null), newE] updates per the diff
$updates = @{
    2 = @("57.156.88", "+1.66%")
    3 = @("3.046.93", "+1.71%")
    4 = @($null, "+0.03%")
    5 = @("514.57", "+2.04%")
    6 = @("139.32", "+0.99%")
    7 = @($null, "-0.03%")
    8 = @($null, "+0.96%")
    9 = @($null, "+2.09%")
    10 = @($null, "+0.77%")
    11 = @($null, "+1.96%")
    12 = @("3.574.61", "+2.08%")
    13 = @($null, "+3.34%")
    14 = @($null, "-3.27%")
    15 = @($null, "+0.73%")
    16 = @("57.282.51", "+1.87%")
    17 = @("6.13", "+1.84%")
    18 = @("3.052.82", "+1.63%")
    19 = @("12.80", "-1.27%")
    20 = @("8.04", "+0.16%")
    21 = @("331.81", "+1.05%")
    22 = @("1.00", "-0.01%")
    23 = @($null, "+0.41%")
    24 = @("65.70", "+1.59%")
    25 = @($null, "+4.38%")
    26 = @($null, "-0.11%")
    27 = @("0.0₃0896", "-0.59%")
    28 = @("6.30", "-2.57%")
    29 = @("7.14", "+1.87%")
    30 = @($null, "+1.73%")
    31 = @("20.73", "+2.54%")
    32 = @("1.17", "+0.46%")
    33 = @("154.69", "+1.16%")
    34 = @("26.96", "+5.33%")
    35 = @("4.41", "-2.78%")
    36 = @("5.88", "+2.08%")
    37 = @($null, "+1.16%")
    38 = @("0.0670", "+1.60%")
    39 = @("3.087.92", "+1.78%")
    40 = @("3.90", "+2.76%")
    41 = @("36.84", "+0.60%")
    42 = @($null, "-0.05%")
    43 = @("0.653", "-0.29%")
    44 = @("2.249.41", "+3.13%")
    45 = @($null, "+7.23%")
    46 = @("1.38", "+1.81%")
    47 = @("20.07", "+1.72%")
    48 = @("5.84", "-0.98%")
    49 = @("0.922", "-0.48%")
    50 = @("260.07", "+13.53%")
    51 = @("0.0875", "+2.73%")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newD = $pair[0]
    $newE = $pair[1]
    if ($null -ne $newD) {
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $newD
    }
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = "  " + $newE + "  "
}